$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the EAC values (column D) for the data rows, keeping the cell styles.
$ws.Range("D3:D169").ClearContents()

# Scroll the view back to the top-left (removes the stored topLeftCell="A2")
# and move the selection to the first data row of column D.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3:D169").Select()
